# "added common utility for runmodes"
# CustomerSuite's Runmode flips from "N" to "Y" (making it run like the
# other suites), which drops "N" from the shared-string table since it's
# no longer referenced anywhere. Also: column A is auto-fit to the longest
# suite name and the active selection ends up on B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CustomerSuite (row 3) Runmode: N -> Y
$ws.Range("B3").Value = "Y"

# Auto-fit column A (SuiteName) to its contents
$ws.Columns.Item(1).AutoFit()

# Leave the cursor on B2, matching the saved selection
$ws.Range("B2").Select()
